$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $addr, $val)
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue $ws 'D2' '62.266.11'
Set-TextValue $ws 'E2' '  -1.40%  '
Set-TextValue $ws 'D3' '2.444.06'
Set-TextValue $ws 'E3' '  -0.18%  '
Set-TextValue $ws 'E4' '  -0.04%  '
Set-TextValue $ws 'D5' '583.03'
Set-TextValue $ws 'E5' '  +1.94%  '
Set-TextValue $ws 'D6' '143.97'
Set-TextValue $ws 'E6' '  -1.74%  '
Set-TextValue $ws 'E8' '  -1.02%  '
Set-TextValue $ws 'D9' '2.441.64'
Set-TextValue $ws 'E9' '  -0.13%  '
Set-TextValue $ws 'E10' '  -3.33%  '
Set-TextValue $ws 'E11' '  +2.55%  '
Set-TextValue $ws 'D12' '5.22'
Set-TextValue $ws 'E12' '  -1.32%  '
Set-TextValue $ws 'E13' '  -2.90%  '
Set-TextValue $ws 'D14' '26.49'
Set-TextValue $ws 'E14' '  -1.66%  '
Set-TextValue $ws 'E15' '  -3.35%  '
Set-TextValue $ws 'D16' '2.881.48'
Set-TextValue $ws 'E16' '  -0.39%  '
Set-TextValue $ws 'D17' '62.144.32'
Set-TextValue $ws 'E17' '  -1.12%  '
Set-TextValue $ws 'D18' '2.442.86'
Set-TextValue $ws 'E18' '  -0.03%  '
Set-TextValue $ws 'D19' '10.91'
Set-TextValue $ws 'E19' '  -3.23%  '
Set-TextValue $ws 'E20' '  -2.49%  '
Set-TextValue $ws 'D21' '329.77'
Set-TextValue $ws 'E22' '  -1.99%  '
Set-TextValue $ws 'E23' '  -3.82%  '
Set-TextValue $ws 'E24' '  -0.53%  '
Set-TextValue $ws 'D25' '65.84'
Set-TextValue $ws 'E25' '  +0.23%  '
Set-TextValue $ws 'D26' '9.41'
Set-TextValue $ws 'E26' '  +5.83%  '
Set-TextValue $ws 'D27' '619.53'
Set-TextValue $ws 'E27' '  +1.28%  '
Set-TextValue $ws 'E28' '  -0.43%  '
Set-TextValue $ws 'E29' '  -6.55%  '
Set-TextValue $ws 'E30' '  -0.21%  '
Set-TextValue $ws 'D31' '1.43'
Set-TextValue $ws 'E31' '  -3.84%  '
Set-TextValue $ws 'D32' '8.01'
Set-TextValue $ws 'E32' '  -2.54%  '
Set-TextValue $ws 'E33' '  +0.30%  '
Set-TextValue $ws 'E34' '  -0.87%  '
Set-TextValue $ws 'E35' '  -4.86%  '
Set-TextValue $ws 'E36' '  +0.20%  '
Set-TextValue $ws 'B37' 'ImmutableX'
Set-TextValue $ws 'C37' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws 'D37' '1.44'
Set-TextValue $ws 'E37' '  -5.60%  '
Set-TextValue $ws 'B38' 'PolygonEcosystemToken'
Set-TextValue $ws 'C38' 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-TextValue $ws 'D38' '0.379'
Set-TextValue $ws 'E38' '  -0.23%  '
Set-TextValue $ws 'D39' '151.70'
Set-TextValue $ws 'E39' '  +2.95%  '
Set-TextValue $ws 'E40' '  -2.25%  '
Set-TextValue $ws 'D41' '5.27'
Set-TextValue $ws 'E41' '  -2.45%  '
Set-TextValue $ws 'E42' '  -1.26%  '
Set-TextValue $ws 'D43' '42.40'
Set-TextValue $ws 'E43' '  +1.01%  '
Set-TextValue $ws 'E44' '  +0.02%  '
Set-TextValue $ws 'E45' '  -4.76%  '
Set-TextValue $ws 'D46' '143.43'
Set-TextValue $ws 'E46' '  -3.42%  '
Set-TextValue $ws 'E47' '  -3.12%  '
Set-TextValue $ws 'D48' '0.0527'
Set-TextValue $ws 'E48' '  -1.01%  '
Set-TextValue $ws 'D49' '0.600'
Set-TextValue $ws 'E49' '  -0.12%  '
Set-TextValue $ws 'D50' '19.56'
Set-TextValue $ws 'E50' '  -7.46%  '
Set-TextValue $ws 'D51' '0.0₆0238'
Set-TextValue $ws 'E51' '  +8.55%  '
